$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet: column B (Prophet Forecast) and column D (yhat_upper) ---

$updates = @(
    @{ Row = 2;  B = $null; D = 221.0849778471946 },
    @{ Row = 3;  B = 159;   D = 201.7484654015377 },
    @{ Row = 4;  B = 153;   D = 195.5400274973452 },
    @{ Row = 5;  B = 156;   D = 199.8839534004867 },
    @{ Row = 6;  B = 156;   D = 201.3808662571929 },
    @{ Row = 7;  B = 145;   D = 185.1395569624067 },
    @{ Row = 8;  B = $null; D = 181.6949562513354 },
    @{ Row = 9;  B = $null; D = 179.616336125509  },
    @{ Row = 10; B = 134;   D = 175.7337833657078 },
    @{ Row = 11; B = $null; D = 167.4459777250811 },
    @{ Row = 12; B = 110;   D = 152.3396702761764 },
    @{ Row = 13; B = 108;   D = 153.6437735335704 },
    @{ Row = 14; B = 123;   D = 165.7719481597741 },
    @{ Row = 15; B = 148;   D = 190.8652103296676 },
    @{ Row = 16; B = 157;   D = 200.4464218693634 },
    @{ Row = 17; B = 140;   D = 185.9286942386037 },
    @{ Row = 18; B = $null; D = 155.2100107058579 },
    @{ Row = 19; B = 93;    D = 133.6299565829612 },
    @{ Row = 20; B = 108;   D = 154.4480246787629 },
    @{ Row = 21; B = $null; D = 176.4309099129671 }
)

foreach ($u in $updates) {
    if ($null -ne $u.B) {
        $wsForecast.Cells.Item($u.Row, 2).Value = $u.B
    }
    $wsForecast.Cells.Item($u.Row, 4).Value = $u.D
}

# --- Summary sheet updates ---
# These cells hold numeric-looking values but are stored as text in the
# workbook, so force text formatting before/while assigning, then restore
# the default "Normal" style so no stray formatting is left behind.
foreach ($cellInfo in @(
        @{ Addr = "B10"; Text = "1221" },
        @{ Addr = "B11"; Text = "646" },
        @{ Addr = "B14"; Text = "93" }
    )) {
    $rng = $wsSummary.Range($cellInfo.Addr)
    $rng.NumberFormat = "@"
    $rng.Value = $cellInfo.Text
    $rng.Style = "Normal"
}
